$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("country_file_name")

# Fill in the "country data iea" column for the newly-added energy data countries
# (Indonesia, Kazakhstan, Turkmenistan now have IEA country files too)
$ws.Range("C9").Value = "ID.csv"
$ws.Range("C10").Value = "KZ.csv"
$ws.Range("C15").Value = "TM.csv"

# Update the active selection to match where the author left off editing
$ws.Range("I31").Select()
